$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Cosmetic cursor/selection moves that happened on a few already-existing
# sheets while the author was editing the workbook (no data changes there,
# just where the cell cursor ended up when the file was saved).
# -----------------------------------------------------------------------

$wsBLF = $wb.Worksheets.Item(19)      # XBR_BR_6pin_BLF
$wsBLF.Activate()
$wsBLF.Range("F40").Select()

$wsPWR10 = $wb.Worksheets.Item(25)    # X2_PWR_10pin_BLZP
$wsPWR10.Activate()
$wsPWR10.Range("E42").Select()

$wsDC = $wb.Worksheets.Item(26)       # X2_320_DC_1778078 (was the active tab)
$wsDC.Activate()
$wsDC.Range("F5").Select()

# -----------------------------------------------------------------------
# New sheet 1: X4_M1_6pin_BLF  (added right after X2_320_DC_1778078)
# -----------------------------------------------------------------------

$ws27 = $wb.Worksheets.Add($null, $wsDC)
$ws27.Name = "X4_M1_6pin_BLF"

$ws27.Range("A1").Value = "Pin č."
$ws27.Range("B1").Value = "Označení"
$ws27.Range("C1").Value = "Popis"
$ws27.Range("D1").Value = "Průřez"

$ws27.Range("A2").Value = 1
$ws27.Range("B2").Value = "-B"
$ws27.Range("C2").Value = "- BRZDA"
$ws27.Range("D2").Value = "0,5 ~ 2,5 mm2"

$ws27.Range("A3").Value = 2
$ws27.Range("B3").Value = "+B"
$ws27.Range("C3").Value = "+ BRZDA"
$ws27.Range("D3").Value = "0,5 ~ 2,5 mm2"

$ws27.Range("A4").Value = 3
$ws27.Range("B4").Value = "PE"
$ws27.Range("C4").Value = "Uzemnění"
$ws27.Range("D4").Value = "0,5 ~ 2,5 mm2"

$ws27.Range("A5").Value = 4
$ws27.Range("B5").Value = "W"
$ws27.Range("C5").Value = "Fáze W"
$ws27.Range("D5").Value = "0,5 ~ 2,5 mm2"

$ws27.Range("A6").Value = 5
$ws27.Range("B6").Value = "V"
$ws27.Range("C6").Value = "Fáze V"
$ws27.Range("D6").Value = "0,5 ~ 2,5 mm2"

$ws27.Range("A7").Value = 6
$ws27.Range("B7").Value = "U"
$ws27.Range("C7").Value = "Fáze U"
$ws27.Range("D7").Value = "0,5 ~ 2,5 mm2"

# leftover empty styled cell (row 18, col C) carried over from the template
# the other pinout sheets were copy/pasted from - copy the exact format
# from a sibling sheet so no new style entries get created.
[void]($wsBLF.Range("C18").Copy($ws27.Range("C18")))
$ws27.Rows.Item(18).RowHeight = $wsBLF.Rows.Item(18).RowHeight

$ws27.Range("E9").Select()

# -----------------------------------------------------------------------
# New sheet 2: X2_PWR_12pin_BLZ (added right after X4_M1_6pin_BLF) - this
# becomes the new active sheet/tab, same as in the source commit.
# -----------------------------------------------------------------------

$ws28 = $wb.Worksheets.Add($null, $ws27)
$ws28.Name = "X2_PWR_12pin_BLZ"

$ws28.Range("A1").Value = "Pin č."
$ws28.Range("B1").Value = "Označení"
$ws28.Range("C1").Value = "Popis"
$ws28.Range("D1").Value = "Průřez"

$ws28.Range("A2").Value = 1
$ws28.Range("B2").Value = "PE"
$ws28.Range("C2").Value = "Uzemnění"
$ws28.Range("D2").Value = "0,2 ~ 4 mm2"

$ws28.Range("A3").Value = 2
$ws28.Range("B3").Value = "L1 "
$ws28.Range("C3").Value = "Fáze 1"
$ws28.Range("D3").Value = "0,2 ~ 4 mm2"

$ws28.Range("A4").Value = 3
$ws28.Range("B4").Value = "L2 "
$ws28.Range("C4").Value = "Fáze 2"
$ws28.Range("D4").Value = "0,2 ~ 4 mm2"

$ws28.Range("A5").Value = 4
$ws28.Range("B5").Value = "L3 "
$ws28.Range("C5").Value = "Fáze 3"
$ws28.Range("D5").Value = "0,2 ~ 4 mm2"

$ws28.Range("A6").Value = 5
$ws28.Range("B6").Value = "RBin "
$ws28.Range("C6").Value = "Brzdný odpor interní"
$ws28.Range("D6").Value = "0,2 ~ 4 mm2"

$ws28.Range("A7").Value = 6
$ws28.Range("B7").Value = "SEL "
$ws28.Range("C7").Value = "Volba brzdného odporu"
$ws28.Range("D7").Value = "0,2 ~ 4 mm2"

$ws28.Range("A8").Value = 7
$ws28.Range("B8").Value = "RBex "
$ws28.Range("C8").Value = "Brzdný odpor externí"
$ws28.Range("D8").Value = "0,2 ~ 4 mm2"

$ws28.Range("A9").Value = 8
$ws28.Range("B9").Value = "+DC "
$ws28.Range("C9").Value = "+ Napájení meziobvodu"
$ws28.Range("D9").Value = "0,2 ~ 4 mm2"

$ws28.Range("A10").Value = 9
$ws28.Range("B10").Value = "-DC "
$ws28.Range("C10").Value = "- Napájení meziobvodu"
$ws28.Range("D10").Value = "0,2 ~ 4 mm2"

$ws28.Range("A11").Value = 10
$ws28.Range("B11").Value = "PE "
$ws28.Range("C11").Value = "Uzemnění"
$ws28.Range("D11").Value = "0,2 ~ 4 mm2"

$ws28.Range("A12").Value = 11
$ws28.Range("B12").Value = "T+ "
$ws28.Range("C12").Value = "+Termistor"
$ws28.Range("D12").Value = "0,2 ~ 4 mm2"

$ws28.Range("A13").Value = 12
$ws28.Range("B13").Value = "T- "
$ws28.Range("C13").Value = "-Termistor"
$ws28.Range("D13").Value = "0,2 ~ 4 mm2"

$ws28.Range("E19").Select()
